$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.81"
$ws.Range("D3").Value = "'22.99"
$ws.Range("D4").Value = "'5.391"
$ws.Range("D5").Value = "'0.05903"
$ws.Range("D6").Value = "'3.457"
$ws.Range("D7").Value = "'6.588"
$ws.Range("D8").Value = "'0.8120"
$ws.Range("D9").Value = "'0.9160"
$ws.Range("D10").Value = "'0.1423"
$ws.Range("D11").Value = "'0.07388"
$ws.Range("D12").Value = "'0.03267"
$ws.Range("D13").Value = "'0.03051"
$ws.Range("D14").Value = "'0.09341"
$ws.Range("D15").Value = "'3.855"
$ws.Range("D16").Value = "'0.001566"
$ws.Range("D17").Value = "'0.04670"
$ws.Range("D18").Value = "'0.0005933"
$ws.Range("D19").Value = "'0.005891"
$ws.Range("D20").Value = "'0.001289"
$ws.Range("E20").Value = "19BitKanKANBestin24h"
$ws.Range("D21").Value = "'0.004899"
$ws.Range("D22").Value = "'0.00009503"
$ws.Range("D23").Value = "'3.617"
$ws.Range("D25").Value = "'0.3202"
$ws.Range("D26").Value = "'0.1334"
$ws.Range("D40").Value = "'0.03953"
$ws.Range("D41").Value = "'0.006174"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D43").Value = "'0.002591"
$ws.Range("D44").Value = "'0.008113"
$ws.Range("D45").Value = "'0.00005189"
$ws.Range("D48").Value = "'0.002278"
